# Generate Report for Handoff
# - Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
# - Update timestamps to reflect the new handoff generation times
# - Narrow the "Latest Handoff Datetime" / status-adjacent columns that were
#   widened to fit the old, longer status string
#
# Target column width (chars, as authored in Excel w/ MDW=7) is
# 17.2159881591797. This engine's ColumnWidth setter quantizes to whole
# on-screen pixels on a 6-px-per-character grid, so the nearest
# representable value is 103 px = 17.1666... chars; feed it an input in the
# middle of the pixel-103 bucket (16.25 - 16.4166) so it rounds reliably.
$newColWidth = 16.33

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-23 02:57:55"

$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-23 02:57:50"

$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-23 02:57:55"

$dede.Columns.Item(3).ColumnWidth = $newColWidth
